# Updates the "ESTADO DE CUENTA" worker/period table (rows 16-42) with the
# new account-statement data (part 1 of the new period batch, 1603-1609) and
# refreshes the totals used by column G.
#
# Strategy: first blank out the C/D/E text cells for the whole table so the
# old shared strings become unused (and get dropped on save), then re-fill
# every row top-to-bottom in the new order. That makes the saved shared
# string table come out ordered by first appearance in the new layout,
# matching how the source workbook was regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=16; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1603"; F=27600; G=781242},
    @{Row=17; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1603"; F=27600; G=781242},
    @{Row=18; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1603"; F=27600; G=781242},
    @{Row=19; C="1047381296"; D="LUISA SALCEDO PADILLA";         E="1604"; F=27600; G=781242},
    @{Row=20; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1604"; F=27600; G=781242},
    @{Row=21; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1604"; F=27600; G=781242},
    @{Row=22; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1604"; F=27600; G=781242},
    @{Row=23; C="1047381296"; D="LUISA SALCEDO PADILLA";         E="1605"; F=27600; G=781242},
    @{Row=24; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1605"; F=27600; G=781242},
    @{Row=25; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1605"; F=27600; G=781242},
    @{Row=26; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1605"; F=27600; G=781242},
    @{Row=27; C="1047381296"; D="LUISA SALCEDO PADILLA";         E="1606"; F=27600; G=781242},
    @{Row=28; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1606"; F=27600; G=781242},
    @{Row=29; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1606"; F=27600; G=781242},
    @{Row=30; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1606"; F=27600; G=781242},
    @{Row=31; C="1047381296"; D="LUISA SALCEDO PADILLA";         E="1607"; F=27578; G=781242},
    @{Row=32; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1607"; F=27578; G=781242},
    @{Row=33; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1607"; F=27578; G=781242},
    @{Row=34; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1607"; F=27578; G=781242},
    @{Row=35; C="1047381296"; D="LUISA SALCEDO PADILLA";         E="1608"; F=27578; G=781242},
    @{Row=36; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1608"; F=27578; G=781242},
    @{Row=37; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1608"; F=27578; G=781242},
    @{Row=38; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1608"; F=27578; G=781242},
    @{Row=39; C="1047381296"; D="LUISA SALCEDO PADILLA";         E="1609"; F=27578; G=781242},
    @{Row=40; C="73185136";   D="LEONARDO GUARDO BENAVIDES";    E="1609"; F=27578; G=781242},
    @{Row=41; C="45530579";   D="YELENA PATRICIA GUERRA CORPAS"; E="1609"; F=27578; G=781242},
    @{Row=42; C="1143351939"; D="JEYFERSON JOSE CORPA AREVALO";  E="1609"; F=27578; G=781242}
)

# Step 1: clear the text columns so stale shared strings fall out of use.
foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
    $ws.Range("E$r").Value = ""
}

# Step 2: write the new data top to bottom, left to right.
foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
}
